$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @{
    2 = @{ "B"=1.02; "C"=1.022667776422923; "D"=1.032581945477635; "E"=0.9926147277508489; "F"=1.039115856583574; "I"=1.031330526671931; "J"=1.027852322274696; "K"=1.035386761362304; "L"=0.9955398523336033; "M"=1.041901960711133; "N"=1.013362350471271 }
    3 = @{ "B"=1.02; "C"=1.02352424478629; "D"=1.033262072007751; "E"=0.9936372048519304; "F"=1.039977413955828; "I"=1.031476498206074; "J"=1.028347421440468; "K"=1.035876034291458; "L"=0.9963617723202692; "M"=1.042573534941571; "N"=1.013525295867261 }
    4 = @{ "B"=1.02; "C"=1.024078554996472; "D"=1.033701806534017; "E"=0.9942998659930995; "F"=1.040535090966167; "I"=1.031569008283956; "J"=1.02866727219305; "K"=1.036191597054664; "L"=0.9968940712668345; "M"=1.043007640725988; "N"=1.0136305568504 }
    5 = @{ "B"=1.02; "C"=1.02431161392079; "D"=1.033886584744381; "E"=0.9945786998346017; "F"=1.040769582560537; "I"=1.031607433417288; "J"=1.028801613859652; "K"=1.036324011800248; "L"=0.997117960005301; "M"=1.043190030096137; "N"=1.013674766146051 }
    6 = @{ "B"=1.02; "C"=1.024350747079497; "D"=1.033917604714873; "E"=0.9946255319796338; "F"=1.040808957271178; "I"=1.031613857807952; "J"=1.028824163136199; "K"=1.036346230245317; "L"=0.9971555583673453; "M"=1.043220647654871; "N"=1.013682186582415 }
    7 = @{ "B"=1.02; "C"=1.024081669037139; "D"=1.033704275890337; "E"=0.9943035907982488; "F"=1.040538224082929; "I"=1.031569523555057; "J"=1.028669067760033; "K"=1.036193367362554; "L"=0.9968970624462087; "M"=1.043010078250747; "N"=1.01363114774435 }
    8 = @{ "B"=1.02; "C"=1.022957198712853; "D"=1.032811869445454; "E"=0.9929600610674301; "F"=1.039406983182445; "I"=1.031380260273309; "J"=1.028019748589288; "K"=1.035552325793444; "L"=0.995817528259106; "M"=1.04212901413597; "N"=1.013417454754821 }
    9 = @{ "B"=1.02; "C"=1.02097669911972; "D"=1.031236708524069; "E"=0.9906006454969559; "F"=1.037415137063201; "I"=1.03103190965279; "J"=1.026871697309547; "K"=1.034414903469186; "L"=0.9939188001724441; "M"=1.040573105381932; "N"=1.013039572592106 }
    10 = @{ "B"=1.02; "C"=1.019657088472008; "D"=1.030184935775408; "E"=0.989033133672735; "F"=1.03608837966399; "I"=1.030789746500634; "J"=1.026103795218602; "K"=1.033651440363727; "L"=0.9926553831429383; "M"=1.039533662497976; "N"=1.012786780961461 }
    11 = @{ "B"=1.02; "C"=1.019085869511543; "D"=1.029729132432412; "E"=0.988355674866747; "F"=1.0355141704471; "I"=1.030682542613625; "J"=1.025770697603464; "K"=1.033319639883312; "L"=0.9921088820399291; "M"=1.039083075988354; "N"=1.012677117662382 }
    12 = @{ "B"=1.02; "C"=1.018873721785012; "D"=1.029559771514906; "E"=0.9881042295826724; "F"=1.035300928122742; "I"=1.030642370716227; "J"=1.025646882598992; "K"=1.033196213017574; "L"=0.9919059725120875; "M"=1.038915634295688; "N"=1.012636353718266 }
    13 = @{ "B"=1.02; "C"=1.018919226886617; "D"=1.029596102480932; "E"=0.9881581567098651; "F"=1.035346667282947; "I"=1.030651003622375; "J"=1.025673445290211; "K"=1.033222696687187; "L"=0.9919494934313052; "M"=1.038951554423579; "N"=1.012645099079459 }
    14 = @{ "B"=1.02; "C"=1.01906833272735; "D"=1.02971513412395; "E"=0.9883348863814464; "F"=1.035496542851801; "I"=1.03067922915996; "J"=1.025760464799884; "K"=1.033309441073073; "L"=0.9920921077337197; "M"=1.039069236702952; "N"=1.012673748714221 }
    15 = @{ "B"=1.02; "C"=1.01916020556024; "D"=1.029788466174261; "E"=0.9884438009545853; "F"=1.035588892088314; "I"=1.030696573278746; "J"=1.025814068799882; "K"=1.03336286315395; "L"=0.9921799884222134; "M"=1.039141734903247; "N"=1.012691396720483 }
    16 = @{ "B"=1.02; "C"=1.019695002313475; "D"=1.030215178102236; "E"=0.9890781214508737; "F"=1.036126494194923; "I"=1.030796811918468; "J"=1.026125889450715; "K"=1.033673435360662; "L"=0.9926916645766087; "M"=1.039563556011707; "N"=1.012794054710701 }
    17 = @{ "B"=1.02; "C"=1.02003051554559; "D"=1.030482742972693; "E"=0.989476357848556; "F"=1.036463795397445; "I"=1.030859061544422; "J"=1.026321328847339; "K"=1.033867924474243; "L"=0.9930127773699352; "M"=1.039828020175554; "N"=1.012858395303802 }
    18 = @{ "B"=1.02; "C"=1.020226232135527; "D"=1.030638772375838; "E"=0.9897087662937556; "F"=1.036660564917639; "I"=1.030895144323732; "J"=1.026435268252477; "K"=1.033981249279659; "L"=0.9932001317071769; "M"=1.039982229223287; "N"=1.012895904486319 }
    19 = @{ "B"=1.02; "C"=1.020292969349011; "D"=1.030691968101817; "E"=0.9897880325774034; "F"=1.036727662842143; "I"=1.030907409197382; "J"=1.026474108929556; "K"=1.034019870156775; "L"=0.9932640239640975; "M"=1.040034802265508; "N"=1.01290869081314 }
    20 = @{ "B"=1.02; "C"=1.019994516342678; "D"=1.030454039574147; "E"=0.9894336180360679; "F"=1.036427603320829; "I"=1.030852406157784; "J"=1.026300365946896; "K"=1.033847069767311; "L"=0.9929783193494215; "M"=1.039799650671896; "N"=1.012851494191427 }
    21 = @{ "B"=1.02; "C"=1.019024423979127; "D"=1.029680083784152; "E"=0.9882828385668249; "F"=1.035452406989412; "I"=1.030670927142098; "J"=1.025734842136154; "K"=1.033283902010566; "L"=0.9920501090198102; "M"=1.039034584230475; "N"=1.01266531293918 }
    22 = @{ "B"=1.02; "C"=1.018414653485722; "D"=1.029193148145298; "E"=0.9875604150241495; "F"=1.034839520889395; "I"=1.030554790224374; "J"=1.025378767877658; "K"=1.032928767230177; "L"=0.9914670000341481; "M"=1.038553130318212; "N"=1.012548079356194 }
    23 = @{ "B"=1.02; "C"=1.01873788831022; "D"=1.029451311621974; "E"=0.9879432794643023; "F"=1.035164398350718; "I"=1.030616549120589; "J"=1.025567577282791; "K"=1.033117129889919; "L"=0.991776070289318; "M"=1.038808398070047; "N"=1.012610243472466 }
    24 = @{ "B"=1.02; "C"=1.02001078277924; "D"=1.03046700951786; "E"=0.9894529299347244; "F"=1.036443956877767; "I"=1.030855414140053; "J"=1.026309838355756; "K"=1.033856493474567; "L"=0.9929938892766442; "M"=1.039812469779204; "N"=1.012854612567638 }
    25 = @{ "B"=1.02; "C"=1.021488583985937; "D"=1.031644226254298; "E"=0.9912096547607049; "F"=1.037929883479147; "I"=1.031123720487791; "J"=1.027168947628348; "K"=1.034709873485362; "L"=0.9944092447426414; "M"=1.040975734154395; "N"=1.013137419208912 }
}

foreach ($row in $rowData.Keys) {
    $cols = $rowData[$row]
    foreach ($col in $cols.Keys) {
        $addr = "$col$row"
        $ws.Range($addr).Value = $cols[$col]
    }
}

Write-Host "Updated $($rowData.Count) rows"
